# Update "想去人数" (want-to-go count) values in column F
# on the "展览" and "全部类型" worksheets, per refreshed site data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value  = 74
    $ws.Range("F5").Value  = 378
    $ws.Range("F6").Value  = 11261
    $ws.Range("F7").Value  = 633
    $ws.Range("F8").Value  = 101
    $ws.Range("F9").Value  = 12
    $ws.Range("F12").Value = 158
    $ws.Range("F13").Value = 19
    $ws.Range("F15").Value = 45
    $ws.Range("F18").Value = 316
    $ws.Range("F19").Value = 1234
    $ws.Range("F21").Value = 889
    $ws.Range("F22").Value = 108
}
